# Update project configuration test workbook:
# - "individualPhysiologyFile" row now points to the new biometrics file
#   instead of the old physiology file, with an updated description.
# - Selection moves to C7 (last active cell before save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B6").Value = "IndividualBiometrics.xlsx"
$ws.Range("C6").Value = "Name of the excel file with individual biometrics information. Must be located in the ""paramsFolder"""

$ws.Activate()
$ws.Range("C7").Select()
